# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp in A1
# - Swap the Cabo Verde / Suazilandia rows' labels (data order on the sheet
#   stays the same, only the two country names trade places)
# - Refresh the numeric COVID columns (B..H) for the countries whose figures
#   moved in this data pull

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Julio de 2020 a las 20:25"

# --- Swap Cabo Verde / Suazilandia labels ------------------------------
# Row 121 currently shows "Cabo Verde", row 122 currently shows "Suazilandia".
# After the update they trade places (the underlying per-row statistics stay
# attached to their original row, only the displayed country name changes).
$ws.Range("A121").Value = "Suazilandia"
$ws.Range("A122").Value = "Cabo Verde"

# --- Helper to update a whole data row (columns B:H) -------------------
# Positional params only -- this host's PowerShell subset does not bind
# named/-Switch style parameters to function scope.
function Set-RowValues {
    param($Row, $Values)
    $cols = @("B","C","D","E","F","G","H")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $v = $Values[$i]
        if ($null -ne $v) {
            $addr = "$($cols[$i])$Row"
            $ws.Range($addr).Value = $v
        }
    }
}

# Estados Unidos
Set-RowValues 4   @(4461071, 27661, 2151853, 2157512, $null, 631, 151706)

# India
Set-RowValues 6   @(1531783, 49280, 988583, 508976, $null, $null, $null)

# Turquia
Set-RowValues 20  @(227982, 963, 211561, 10776, $null, 15, 5645)

# Francia
Set-RowValues 22  @(183804, 725, $null, 72499, $null, 14, 30223)

# Canada
Set-RowValues 25  @(114877, 280, 99968, 6001, $null, 7, 8908)

# Kazajistan
Set-RowValues 29  @($null, $null, $null, 29451, $null, $null, 793)

# Argelia
Set-RowValues 60  @($null, $null, 19233, 8208, $null, $null, $null)

# Moldavia
Set-RowValues 63  @(23521, 367, $null, 6306, $null, 5, 753)

# Marruecos
Set-RowValues 65  @(21387, 500, 17066, 3994, $null, 11, 327)

# Guayana Francesa
Set-RowValues 89  @(7562, 48, 6106, 1413, $null, 1, 43)

# Republica de Yibuti
Set-RowValues 97  @(5068, 9, 4992, 18, $null, $null, $null)

# Sri Lanka
Set-RowValues 117 @(2809, 4, $null, 502, $null, $null, $null)

# Row 121 (now "Suazilandia" label, same row/statistics as before)
Set-RowValues 121 @(2404, 88, 1025, 1340, $null, 5, 39)

# Row 122 (now "Cabo Verde" label, same row/statistics as before)
Set-RowValues 122 @(2354, 26, 1616, 716, $null, $null, 22)

# Republica de Chipre
Set-RowValues 147 @(1067, 7, $null, 196, $null, $null, $null)

# Monaco
Set-RowValues 184 @(117, 1, $null, 9, $null, $null, $null)
